$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.197.04"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.271.74"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.24%  "
$ws.Range("E7").Value = "  -2.31%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.95%  "
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.69%  "
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "2.624.77"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").Value = "2.268.57"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.785"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("D19").Value = "42.153.30"
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "0.0₃0892"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -2.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("E30").Value = "  -4.82%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.34%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +5.18%  "
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("E37").Value = "  -3.16%  "
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0985"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("E42").Value = "  -4.36%  "
$ws.Range("E43").Value = "  -7.06%  "
$ws.Range("D44").Value = "1.959.72"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.21%  "
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("D49").Value = "2.495.94"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.53%  "
$ws.Range("E51").Value = "  -3.53%  "
